$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values of columns A, Q, R between row 2 and row 3
$colsToSwap = @("A", "Q", "R")

foreach ($col in $colsToSwap) {
    $addr2 = "${col}2"
    $addr3 = "${col}3"
    $cell2 = $ws.Range($addr2)
    $cell3 = $ws.Range($addr3)
    $temp = $cell2.Value2
    $cell2.Value2 = $cell3.Value2
    $cell3.Value2 = $temp
}
